$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped from 45175 to 45177
# for every data row (rows 2 through 300).
$ws.Range("C2:C300").Value = 45177
